# Refactor synthetic array /3 for publipostage
# Mapping of old symbols/labels to new ones:
#   black square  -> blue book, "noir" -> "bleu"
#   red square    -> red book   (label "rouge" unchanged)
#   green square  -> green book (label "vert" unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blackSquare = "⬛"
$redSquare   = "🟥"
$greenSquare = "🟩"

$blueBook  = "📘"
$redBook   = "📕"
$greenBook = "📗"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$firstRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2

    if ($valA -eq $blackSquare) {
        $cellA.Value = $blueBook
        if ($cellB.Value2 -eq "noir") {
            $cellB.Value = "bleu"
        }
    }
    elseif ($valA -eq $redSquare) {
        $cellA.Value = $redBook
    }
    elseif ($valA -eq $greenSquare) {
        $cellA.Value = $greenBook
    }
}
